$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 7).Value2 = 1.625602666666667
$ws.Cells.Item(2, 8).Value2 = 4.876808
$ws.Cells.Item(2, 9).Value2 = 0.08794732633208746
$ws.Cells.Item(2, 10).Value2 = 0.08794732633208747
$ws.Cells.Item(2, 13).Value2 = 0.02507166666666667
$ws.Cells.Item(2, 14).Value2 = 0.075215
$ws.Cells.Item(2, 15).Value2 = 0.009392568139045224
$ws.Cells.Item(2, 16).Value2 = 0.009392568139045224
$ws.Cells.Item(2, 17).Value2 = 0.04075656819111112
$ws.Cells.Item(2, 18).Value2 = 0.3668091137200001
$ws.Cells.Item(2, 19).Value2 = 0.0008260512552209778
$ws.Cells.Item(2, 20).Value2 = 0.0008260512552209779
$ws.Cells.Item(3, 7).Value2 = 1.625602666666667
$ws.Cells.Item(3, 8).Value2 = 4.876808
$ws.Cells.Item(3, 9).Value2 = 0.08794732633208746
$ws.Cells.Item(3, 10).Value2 = 0.08794732633208747
$ws.Cells.Item(3, 14).Value2 = 7.038411000000001
$ws.Cells.Item(3, 15).Value2 = 0.8789304647757153
$ws.Cells.Item(3, 16).Value2 = 0.8789304647757155
$ws.Cells.Item(3, 17).Value2 = 3.813886563565334
$ws.Cells.Item(3, 18).Value2 = 34.32497907208801
$ws.Cells.Item(3, 19).Value2 = 0.07729958440884314
$ws.Cells.Item(3, 20).Value2 = 0.07729958440884316
$ws.Cells.Item(4, 7).Value2 = 1.625602666666667
$ws.Cells.Item(4, 8).Value2 = 4.876808
$ws.Cells.Item(4, 9).Value2 = 0.08794732633208746
$ws.Cells.Item(4, 10).Value2 = 0.08794732633208747
$ws.Cells.Item(4, 13).Value2 = 0.2981003333333334
$ws.Cells.Item(4, 14).Value2 = 0.894301
$ws.Cells.Item(4, 15).Value2 = 0.1116769670852394
$ws.Cells.Item(4, 16).Value2 = 0.1116769670852394
$ws.Cells.Item(4, 17).Value2 = 0.484592696800889
$ws.Cells.Item(4, 18).Value2 = 4.361334271208
$ws.Cells.Item(4, 19).Value2 = 0.009821690668023341
$ws.Cells.Item(4, 20).Value2 = 0.009821690668023342
$ws.Cells.Item(5, 9).Value2 = 0.2830836711542908
$ws.Cells.Item(5, 10).Value2 = 0.2830836711542908
$ws.Cells.Item(5, 13).Value2 = 0.02507166666666667
$ws.Cells.Item(5, 14).Value2 = 0.075215
$ws.Cells.Item(5, 15).Value2 = 0.009392568139045224
$ws.Cells.Item(5, 16).Value2 = 0.009392568139045224
$ws.Cells.Item(5, 17).Value2 = 0.13118669354
$ws.Cells.Item(5, 18).Value2 = 1.18068024186
$ws.Cells.Item(5, 19).Value2 = 0.002658882670367747
$ws.Cells.Item(5, 20).Value2 = 0.002658882670367747
$ws.Cells.Item(6, 9).Value2 = 0.2830836711542908
$ws.Cells.Item(6, 10).Value2 = 0.2830836711542908
$ws.Cells.Item(6, 14).Value2 = 7.038411000000001
$ws.Cells.Item(6, 15).Value2 = 0.8789304647757153
$ws.Cells.Item(6, 16).Value2 = 0.8789304647757155
$ws.Cells.Item(6, 19).Value2 = 0.2488108626580566
$ws.Cells.Item(6, 20).Value2 = 0.2488108626580566
$ws.Cells.Item(7, 9).Value2 = 0.2830836711542908
$ws.Cells.Item(7, 10).Value2 = 0.2830836711542908
$ws.Cells.Item(7, 13).Value2 = 0.2981003333333334
$ws.Cells.Item(7, 14).Value2 = 0.894301
$ws.Cells.Item(7, 15).Value2 = 0.1116769670852394
$ws.Cells.Item(7, 16).Value2 = 0.1116769670852394
$ws.Cells.Item(7, 17).Value2 = 1.559800454956
$ws.Cells.Item(7, 18).Value2 = 14.038204094604
$ws.Cells.Item(7, 19).Value2 = 0.03161392582586647
$ws.Cells.Item(7, 20).Value2 = 0.03161392582586647
$ws.Cells.Item(8, 7).Value2 = 3.405616666666667
$ws.Cells.Item(8, 8).Value2 = 10.21685
$ws.Cells.Item(8, 9).Value2 = 0.1842485168651273
$ws.Cells.Item(8, 10).Value2 = 0.1842485168651273
$ws.Cells.Item(8, 13).Value2 = 0.02507166666666667
$ws.Cells.Item(8, 14).Value2 = 0.075215
$ws.Cells.Item(8, 15).Value2 = 0.009392568139045224
$ws.Cells.Item(8, 16).Value2 = 0.009392568139045224
$ws.Cells.Item(8, 17).Value2 = 0.08538448586111111
$ws.Cells.Item(8, 18).Value2 = 0.76846037275
$ws.Cells.Item(8, 19).Value2 = 0.001730566749173731
$ws.Cells.Item(8, 20).Value2 = 0.001730566749173731
$ws.Cells.Item(9, 7).Value2 = 3.405616666666667
$ws.Cells.Item(9, 8).Value2 = 10.21685
$ws.Cells.Item(9, 9).Value2 = 0.1842485168651273
$ws.Cells.Item(9, 10).Value2 = 0.1842485168651273
$ws.Cells.Item(9, 14).Value2 = 7.038411000000001
$ws.Cells.Item(9, 15).Value2 = 0.8789304647757153
$ws.Cells.Item(9, 16).Value2 = 0.8789304647757155
$ws.Cells.Item(9, 17).Value2 = 7.990043269483333
$ws.Cells.Item(9, 18).Value2 = 71.91038942535
$ws.Cells.Item(9, 19).Value2 = 0.1619416345625025
$ws.Cells.Item(9, 20).Value2 = 0.1619416345625026
$ws.Cells.Item(10, 7).Value2 = 3.405616666666667
$ws.Cells.Item(10, 8).Value2 = 10.21685
$ws.Cells.Item(10, 9).Value2 = 0.1842485168651273
$ws.Cells.Item(10, 10).Value2 = 0.1842485168651273
$ws.Cells.Item(10, 13).Value2 = 0.2981003333333334
$ws.Cells.Item(10, 14).Value2 = 0.894301
$ws.Cells.Item(10, 15).Value2 = 0.1116769670852394
$ws.Cells.Item(10, 16).Value2 = 0.1116769670852394
$ws.Cells.Item(10, 17).Value2 = 1.015215463538889
$ws.Cells.Item(10, 18).Value2 = 9.136939171849999
$ws.Cells.Item(10, 19).Value2 = 0.020576315553451
$ws.Cells.Item(10, 20).Value2 = 0.020576315553451
$ws.Cells.Item(11, 7).Value2 = 5.048711333333333
$ws.Cells.Item(11, 8).Value2 = 15.146134
$ws.Cells.Item(11, 9).Value2 = 0.273142184307343
$ws.Cells.Item(11, 10).Value2 = 0.273142184307343
$ws.Cells.Item(11, 13).Value2 = 0.02507166666666667
$ws.Cells.Item(11, 14).Value2 = 0.075215
$ws.Cells.Item(11, 15).Value2 = 0.009392568139045224
$ws.Cells.Item(11, 16).Value2 = 0.009392568139045224
$ws.Cells.Item(11, 17).Value2 = 0.1265796076455556
$ws.Cells.Item(11, 18).Value2 = 1.13921646881
$ws.Cells.Item(11, 19).Value2 = 0.002565506577754369
$ws.Cells.Item(11, 20).Value2 = 0.002565506577754369
$ws.Cells.Item(12, 7).Value2 = 5.048711333333333
$ws.Cells.Item(12, 8).Value2 = 15.146134
$ws.Cells.Item(12, 9).Value2 = 0.273142184307343
$ws.Cells.Item(12, 10).Value2 = 0.273142184307343
$ws.Cells.Item(12, 14).Value2 = 7.038411000000001
$ws.Cells.Item(12, 15).Value2 = 0.8789304647757153
$ws.Cells.Item(12, 16).Value2 = 0.8789304647757155
$ws.Cells.Item(12, 17).Value2 = 11.84496846145267
$ws.Cells.Item(12, 18).Value2 = 106.604716153074
$ws.Cells.Item(12, 19).Value2 = 0.2400729870031071
$ws.Cells.Item(12, 20).Value2 = 0.2400729870031072
$ws.Cells.Item(13, 7).Value2 = 5.048711333333333
$ws.Cells.Item(13, 8).Value2 = 15.146134
$ws.Cells.Item(13, 9).Value2 = 0.273142184307343
$ws.Cells.Item(13, 10).Value2 = 0.273142184307343
$ws.Cells.Item(13, 13).Value2 = 0.2981003333333334
$ws.Cells.Item(13, 14).Value2 = 0.894301
$ws.Cells.Item(13, 15).Value2 = 0.1116769670852394
$ws.Cells.Item(13, 16).Value2 = 0.1116769670852394
$ws.Cells.Item(13, 17).Value2 = 1.505022531370445
$ws.Cells.Item(13, 18).Value2 = 13.545202782334
$ws.Cells.Item(13, 19).Value2 = 0.03050369072648155
$ws.Cells.Item(13, 20).Value2 = 0.03050369072648155
$ws.Cells.Item(14, 5).Value2 = 3
$ws.Cells.Item(14, 6).Value2 = 1
$ws.Cells.Item(14, 7).Value2 = 2.689098
$ws.Cells.Item(14, 8).Value2 = 8.067294
$ws.Cells.Item(14, 9).Value2 = 0.145483877576253
$ws.Cells.Item(14, 10).Value2 = 0.145483877576253
$ws.Cells.Item(14, 13).Value2 = 0.02507166666666667
$ws.Cells.Item(14, 14).Value2 = 0.075215
$ws.Cells.Item(14, 15).Value2 = 0.009392568139045224
$ws.Cells.Item(14, 16).Value2 = 0.009392568139045224
$ws.Cells.Item(14, 17).Value2 = 0.06742016869
$ws.Cells.Item(14, 18).Value2 = 0.60678151821
$ws.Cells.Item(14, 19).Value2 = 0.001366467233267469
$ws.Cells.Item(14, 20).Value2 = 0.00136646723326747
$ws.Cells.Item(15, 5).Value2 = 3
$ws.Cells.Item(15, 6).Value2 = 1
$ws.Cells.Item(15, 7).Value2 = 2.689098
$ws.Cells.Item(15, 8).Value2 = 8.067294
$ws.Cells.Item(15, 9).Value2 = 0.145483877576253
$ws.Cells.Item(15, 10).Value2 = 0.145483877576253
$ws.Cells.Item(15, 14).Value2 = 7.038411000000001
$ws.Cells.Item(15, 15).Value2 = 0.8789304647757153
$ws.Cells.Item(15, 16).Value2 = 0.8789304647757155
$ws.Cells.Item(15, 17).Value2 = 6.308992314426001
$ws.Cells.Item(15, 18).Value2 = 56.78093082983401
$ws.Cells.Item(15, 19).Value2 = 0.1278702121354693
$ws.Cells.Item(15, 20).Value2 = 0.1278702121354693
$ws.Cells.Item(16, 5).Value2 = 3
$ws.Cells.Item(16, 6).Value2 = 1
$ws.Cells.Item(16, 7).Value2 = 2.689098
$ws.Cells.Item(16, 8).Value2 = 8.067294
$ws.Cells.Item(16, 9).Value2 = 0.145483877576253
$ws.Cells.Item(16, 10).Value2 = 0.145483877576253
$ws.Cells.Item(16, 13).Value2 = 0.2981003333333334
$ws.Cells.Item(16, 14).Value2 = 0.894301
$ws.Cells.Item(16, 15).Value2 = 0.1116769670852394
$ws.Cells.Item(16, 16).Value2 = 0.1116769670852394
$ws.Cells.Item(16, 17).Value2 = 0.8016210101660001
$ws.Cells.Item(16, 18).Value2 = 7.214589091494
$ws.Cells.Item(16, 19).Value2 = 0.0162471982075162
$ws.Cells.Item(16, 20).Value2 = 0.01624719820751621
$ws.Cells.Item(17, 5).Value2 = 2
$ws.Cells.Item(17, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(17, 7).Value2 = 0.4823246666666667
$ws.Cells.Item(17, 8).Value2 = 1.446974
$ws.Cells.Item(17, 9).Value2 = 0.02609442376489825
$ws.Cells.Item(17, 10).Value2 = 0.02609442376489825
$ws.Cells.Item(17, 13).Value2 = 0.02507166666666667
$ws.Cells.Item(17, 14).Value2 = 0.075215
$ws.Cells.Item(17, 15).Value2 = 0.009392568139045224
$ws.Cells.Item(17, 16).Value2 = 0.009392568139045224
$ws.Cells.Item(17, 17).Value2 = 0.01209268326777778
$ws.Cells.Item(17, 18).Value2 = 0.10883414941
$ws.Cells.Item(17, 19).Value2 = 0.0002450936532609278
$ws.Cells.Item(17, 20).Value2 = 0.0002450936532609279
$ws.Cells.Item(18, 5).Value2 = 2
$ws.Cells.Item(18, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(18, 7).Value2 = 0.4823246666666667
$ws.Cells.Item(18, 8).Value2 = 1.446974
$ws.Cells.Item(18, 9).Value2 = 0.02609442376489825
$ws.Cells.Item(18, 10).Value2 = 0.02609442376489825
$ws.Cells.Item(18, 14).Value2 = 7.038411000000001
$ws.Cells.Item(18, 15).Value2 = 0.8789304647757153
$ws.Cells.Item(18, 16).Value2 = 0.8789304647757155
$ws.Cells.Item(18, 17).Value2 = 1.131599746479333
$ws.Cells.Item(18, 18).Value2 = 10.184397718314
$ws.Cells.Item(18, 19).Value2 = 0.02293518400773649
$ws.Cells.Item(18, 20).Value2 = 0.0229351840077365
$ws.Cells.Item(19, 5).Value2 = 2
$ws.Cells.Item(19, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(19, 7).Value2 = 0.4823246666666667
$ws.Cells.Item(19, 8).Value2 = 1.446974
$ws.Cells.Item(19, 9).Value2 = 0.02609442376489825
$ws.Cells.Item(19, 10).Value2 = 0.02609442376489825
$ws.Cells.Item(19, 13).Value2 = 0.2981003333333334
$ws.Cells.Item(19, 14).Value2 = 0.894301
$ws.Cells.Item(19, 15).Value2 = 0.1116769670852394
$ws.Cells.Item(19, 16).Value2 = 0.1116769670852394
$ws.Cells.Item(19, 17).Value2 = 0.1437811439082222
$ws.Cells.Item(19, 18).Value2 = 1.294030295174
$ws.Cells.Item(19, 19).Value2 = 0.002914146103900831
$ws.Cells.Item(19, 20).Value2 = 0.002914146103900831
